$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the single remaining data row (row 3) with the new movement ---
# Values are entered as a text formula + paste-values so they land as literal
# text (matching the report's existing "amount/date as text" convention)
# instead of being auto-converted to currency/date numbers.

$ws.Range("A3").Formula = "=""LIBRO Y VIAJE """
$ws.Range("A3").Copy()
$ws.Range("A3").PasteSpecial(-4163)

$ws.Range("B3").Formula = "=""16-02-2018"""
$ws.Range("B3").Copy()
$ws.Range("B3").PasteSpecial(-4163)

$ws.Range("D3").Formula = "=""$ 2500.00"""
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)

# C3 ("INGRESO") is unchanged.

# --- Remove all the other movement rows (4-11); only one entry remains now ---
$ws.Rows("4:11").Delete()

# --- Column A's width was sized to fit its contents; recalculate now that
#     the longest value in the column changed ---
$ws.Columns("A:A").EntireColumn.AutoFit()

Write-Output "done"
